# Updated cryptos list on Tue May 14 17:58:25 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.263.13"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.875.48"
$ws.Range("E3").Value = "  -2.25%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.54"
$ws.Range("E5").Value = "  -4.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.04"
$ws.Range("E6").Value = "  -2.80%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.874.96"
$ws.Range("E9").Value = "  -2.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.84"
$ws.Range("E10").Value = "  -6.60%  "
$ws.Range("E11").Value = "  -3.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.431"
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("E13").Value = "  -2.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.79"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.347.93"
$ws.Range("E16").Value = "  -2.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.243.78"
$ws.Range("E17").Value = "  -2.39%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.58"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.889.54"
$ws.Range("E19").Value = "  -2.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "430.76"
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.03"
$ws.Range("E21").Value = "  -2.83%  "
$ws.Range("E22").Value = "  -2.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.81"
$ws.Range("E23").Value = "  -3.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.03"
$ws.Range("E24").Value = "  -2.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.72"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.95"
$ws.Range("E27").Value = "  -10.45%  "
$ws.Range("E28").Value = "  -6.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000102"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.00"
$ws.Range("E30").Value = "  -3.29%  "
$ws.Range("E31").Value = "  -4.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("E32").Value = "  -9.01%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("E34").Value = "  -2.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.46"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.953"
$ws.Range("E36").Value = "  -3.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.37"
$ws.Range("E37").Value = "  -3.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.80"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  -10.59%  "
$ws.Range("E41").Value = "  -3.35%  "
$ws.Range("E42").Value = "  -3.41%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "39.22"
$ws.Range("E43").Value = "  -0.54%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.267"
$ws.Range("E44").Value = "  -4.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.686.27"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.63"
$ws.Range("E46").Value = "  -0.74%  "
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "338.55"
$ws.Range("E49").Value = "  -6.69%  "
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("E51").Value = "  -6.13%  "
